$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Summary_A")

# Insert two new columns at C:D (shifts old C..H to E..J)
$ws.Range("C:D").Insert()

# Copy style of existing data/header columns into the new K:M range so they inherit s=9 / s=10
$ws.Range("J1:J41").Copy()
$ws.Range("K1:M41").PasteSpecial(-4122)

# Set new column widths (C:D and I:M) to match the existing data columns
$ws.Range("C1:D1").ColumnWidth = 13.17
$ws.Range("I1:M1").ColumnWidth = 13.17

# Header row
$ws.Range("C1").Value = "t1/2"
$ws.Range("D1").Value = "unit"
$ws.Range("K1").Value = "A(Bq)@t=15 m"
$ws.Range("L1").Value = "A(Bq)@t=1 hour"
$ws.Range("M1").Value = "A(Bq)@t=1 day"

# Data rows: t1/2 (C), unit (D), and new time points K, L, M
$ws.Range("C2").Value = 58
$ws.Range("D2").Value = "m"
$ws.Range("K2").Value = 117.83
$ws.Range("L2").Value = 68.83
$ws.Range("M2").Value = 0
$ws.Range("C3").Value = 39.6
$ws.Range("D3").Value = "m"
$ws.Range("K3").Value = 137.71
$ws.Range("L3").Value = 62.66
$ws.Range("M3").Value = 0
$ws.Range("C4").Value = 44.3
$ws.Range("D4").Value = "s"
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("C5").Value = 6.5
$ws.Range("D5").Value = "h"
$ws.Range("K5").Value = 407.62
$ws.Range("L5").Value = 376.29
$ws.Range("M5").Value = 32.4
$ws.Range("C6").Value = 32.4
$ws.Range("D6").Value = "m"
$ws.Range("K6").Value = 4074.14
$ws.Range("L6").Value = 1556.06
$ws.Range("M6").Value = 0
$ws.Range("C7").Value = 50.4
$ws.Range("D7").Value = "s"
$ws.Range("K7").Value = 0.03
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("C8").Value = 23.96
$ws.Range("D8").Value = "m"
$ws.Range("K8").Value = 920.82
$ws.Range("L8").Value = 250.56
$ws.Range("M8").Value = 0
$ws.Range("C9").Value = 8.28
$ws.Range("D9").Value = "d"
$ws.Range("K9").Value = 5.99
$ws.Range("L9").Value = 5.98
$ws.Range("M9").Value = 5.52
$ws.Range("C10").Value = 6.2
$ws.Range("D10").Value = "m"
$ws.Range("K10").Value = 1793.05
$ws.Range("L10").Value = 11.73
$ws.Range("M10").Value = 0
$ws.Range("C11").Value = 5.2
$ws.Range("D11").Value = "m"
$ws.Range("K11").Value = 1300.02
$ws.Range("L11").Value = 3.23
$ws.Range("M11").Value = 0
$ws.Range("C12").Value = 41.29
$ws.Range("D12").Value = "d"
$ws.Range("K12").Value = 1.61
$ws.Range("L12").Value = 1.61
$ws.Range("M12").Value = 1.58
$ws.Range("C13").Value = 7.23
$ws.Range("D13").Value = "m"
$ws.Range("K13").Value = 547.97
$ws.Range("L13").Value = 7.34
$ws.Range("M13").Value = 0
$ws.Range("C14").Value = 55.5
$ws.Range("D14").Value = "m"
$ws.Range("K14").Value = 4391.54
$ws.Range("L14").Value = 2503.75
$ws.Range("M14").Value = 0
$ws.Range("C15").Value = 5.07
$ws.Range("D15").Value = "m"
$ws.Range("K15").Value = 352.39
$ws.Range("L15").Value = 0.75
$ws.Range("M15").Value = 0
$ws.Range("C16").Value = 48
$ws.Range("D16").Value = "s"
$ws.Range("K16").Value = 0.01
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 0
$ws.Range("C17").Value = 69.2
$ws.Range("D17").Value = "m"
$ws.Range("K17").Value = 1791.78
$ws.Range("L17").Value = 1141.75
$ws.Range("M17").Value = 0
$ws.Range("C18").Value = 33.5
$ws.Range("D18").Value = "m"
$ws.Range("K18").Value = 2403.21
$ws.Range("L18").Value = 947.35
$ws.Range("M18").Value = 0
$ws.Range("C19").Value = 57.7
$ws.Range("D19").Value = "m"
$ws.Range("K19").Value = 401.03
$ws.Range("L19").Value = 233.59
$ws.Range("M19").Value = 0
$ws.Range("C20").Value = 1.8
$ws.Range("D20").Value = "m"
$ws.Range("K20").Value = 0.28
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("C21").Value = 65.7
$ws.Range("D21").Value = "m"
$ws.Range("K21").Value = 3485.98
$ws.Range("L21").Value = 2168.62
$ws.Range("M21").Value = 0
$ws.Range("C22").Value = 16.991
$ws.Range("D22").Value = "d"
$ws.Range("K22").Value = 5.81
$ws.Range("L22").Value = 5.8
$ws.Range("M22").Value = 5.58
$ws.Range("C23").Value = 56.114
$ws.Range("D23").Value = "m"
$ws.Range("K23").Value = 80.61
$ws.Range("L23").Value = 46.24
$ws.Range("M23").Value = 0
$ws.Range("C24").Value = 39.247
$ws.Range("D24").Value = "d"
$ws.Range("K24").Value = 0.02
$ws.Range("L24").Value = 0.02
$ws.Range("M24").Value = 0.02
$ws.Range("C25").Value = 12.9
$ws.Range("D25").Value = "m"
$ws.Range("K25").Value = 1962.67
$ws.Range("L25").Value = 174.97
$ws.Range("M25").Value = 0
$ws.Range("C26").Value = 7.7
$ws.Range("D26").Value = "m"
$ws.Range("K26").Value = 1180.85
$ws.Range("L26").Value = 20.57
$ws.Range("M26").Value = 0
$ws.Range("C27").Value = 207.3
$ws.Range("D27").Value = "d"
$ws.Range("K27").Value = 0.04
$ws.Range("L27").Value = 0.04
$ws.Range("M27").Value = 0.04
$ws.Range("C28").Value = 3.742
$ws.Range("D28").Value = "y"
$ws.Range("K28").Value = 0.01
$ws.Range("L28").Value = 0.01
$ws.Range("M28").Value = 0.01
$ws.Range("C29").Value = 11.1
$ws.Range("D29").Value = "m"
$ws.Range("K29").Value = 34.06
$ws.Range("L29").Value = 2.05
$ws.Range("M29").Value = 0
$ws.Range("C30").Value = 8.47
$ws.Range("D30").Value = "h"
$ws.Range("K30").Value = 15.45
$ws.Range("L30").Value = 14.53
$ws.Range("M30").Value = 2.21
$ws.Range("C31").Value = 4.34
$ws.Range("D31").Value = "d"
$ws.Range("K31").Value = 3.53
$ws.Range("L31").Value = 3.52
$ws.Range("M31").Value = 3.02
$ws.Range("C32").Value = 20.5
$ws.Range("D32").Value = "h"
$ws.Range("K32").Value = 59.77
$ws.Range("L32").Value = 58.28
$ws.Range("M32").Value = 26.78
$ws.Range("C33").Value = 4.6
$ws.Range("D33").Value = "m"
$ws.Range("K33").Value = 189.27
$ws.Range("L33").Value = 0.22
$ws.Range("M33").Value = 0
$ws.Range("C34").Value = 16.1
$ws.Range("D34").Value = "d"
$ws.Range("K34").Value = 2.14
$ws.Range("L34").Value = 2.14
$ws.Range("M34").Value = 2.05
$ws.Range("C35").Value = 4.7
$ws.Range("D35").Value = "h"
$ws.Range("K35").Value = 158.16
$ws.Range("L35").Value = 141.6
$ws.Range("M35").Value = 4.77
$ws.Range("C36").Value = 8.72
$ws.Range("D36").Value = "m"
$ws.Range("K36").Value = 6.7
$ws.Range("L36").Value = 0.19
$ws.Range("M36").Value = 0
$ws.Range("C37").Value = 3.6
$ws.Range("D37").Value = "m"
$ws.Range("K37").Value = 1.24
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("C38").Value = 91
$ws.Range("D38").Value = "d"
$ws.Range("K38").Value = 0.03
$ws.Range("L38").Value = 0.03
$ws.Range("M38").Value = 0.03
$ws.Range("C39").Value = 4.28
$ws.Range("D39").Value = "d"
$ws.Range("K39").Value = 0.3
$ws.Range("L39").Value = 0.3
$ws.Range("M39").Value = 0.26
$ws.Range("C40").Value = 51.5
$ws.Range("D40").Value = "m"
$ws.Range("K40").Value = 20.14
$ws.Range("L40").Value = 10.99
$ws.Range("M40").Value = 0
$ws.Range("C41").Value = 6.263
$ws.Range("D41").Value = "m"
$ws.Range("K41").Value = 3.85
$ws.Range("L41").Value = 0.03
$ws.Range("M41").Value = 0

# Rename sheet
$ws.Name = "Summary_Act"
